$wb = $excel.ActiveWorkbook

# --- Overview sheet: row 3 is the d4543783-... file ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E3").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F3").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: row 3 is the d4543783-... file ---
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C3").Value = "Handed back: in sync with en-US"
$wsZh.Range("K3").Value = "2016-08-27 00:48:46"
$wsZh.Range("P3").Value = ""

# --- de-de sheet: row 3 is the d4543783-... file ---
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C3").Value = "Handed back: in sync with en-US"
$wsDe.Range("K3").Value = "2016-08-27 00:48:53"
$wsDe.Range("P3").Value = ""
